$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 13 new rows at row 30 (shifts nothing below since sheet ends at row 29)
$ws.Rows("30:42").Insert()

# Row 30
$ws.Range("A30").Value = '0.3333333333333333-FlowTest150424nanana'
$ws.Range("B30").Value = 'Failed'
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = '09:29:30'
$ws.Range("AA1").Formula = "=""2024-04-16"""
$ws.Range("AA1").Copy()
$ws.Range("E30").PasteSpecial(-4163)
$ws.Range("F30").Value = '[1, 2]'
$ws.Range("G30").Value = 0.3333333333333333
$ws.Range("H30").Value = 120
$ws.Range("I30").Value = 50
$ws.Range("J30").Value = 29.00098133087158
$ws.Range("K30").Value = 'BufferX'
$ws.Range("L30").Value = 50
$ws.Range("M30").Value = 2.124137931034483
$ws.Range("N30").Value = 'Lipid1'
$ws.Range("O30").Value = 'na'
$ws.Range("P30").Value = 50
$ws.Range("Q30").Value = 1.26071
$ws.Range("R30").Value = 'Lipid2'
$ws.Range("S30").Value = 'na'
$ws.Range("T30").Value = 50
$ws.Range("U30").Value = 0.9339168965517243
$ws.Range("V30").Value = 'Lipid3'
$ws.Range("W30").Value = 'na'
$ws.Range("X30").Value = 50
$ws.Range("Y30").Value = 1.042847931034482

# Row 31
$ws.Range("A31").Value = '0.3333333333333333-FlowTest150424nanana'
$ws.Range("B31").Value = 'Failed'
$ws.Range("C31").Value = 1
$ws.Range("D31").Value = '10:11:25'
$ws.Range("AA1").Formula = "=""2024-04-16"""
$ws.Range("AA1").Copy()
$ws.Range("E31").PasteSpecial(-4163)
$ws.Range("F31").Value = '[1, 2]'
$ws.Range("G31").Value = 0.3333333333333333
$ws.Range("H31").Value = 120
$ws.Range("I31").Value = 50
$ws.Range("J31").Value = 64.67762470245361
$ws.Range("K31").Value = 'BufferX'
$ws.Range("L31").Value = 50
$ws.Range("M31").Value = 1.256551724137931
$ws.Range("N31").Value = 'Lipid1'
$ws.Range("O31").Value = 'na'
$ws.Range("P31").Value = 50
$ws.Range("Q31").Value = 1.101503103448275
$ws.Range("R31").Value = 'Lipid2'
$ws.Range("S31").Value = 'na'
$ws.Range("T31").Value = 50
$ws.Range("U31").Value = 0.8361582758620688
$ws.Range("V31").Value = 'Lipid3'
$ws.Range("W31").Value = 'na'
$ws.Range("X31").Value = 50
$ws.Range("Y31").Value = 0.763537586206896

# Row 32
$ws.Range("A32").Value = '0.3333333333333333-FlowTest150424nanana'
$ws.Range("B32").Value = 'Failed'
$ws.Range("C32").Value = 1
$ws.Range("D32").Value = '10:12:21'
$ws.Range("AA1").Formula = "=""2024-04-16"""
$ws.Range("AA1").Copy()
$ws.Range("E32").PasteSpecial(-4163)
$ws.Range("F32").Value = '[1, 3]'
$ws.Range("G32").Value = 0.3333333333333333
$ws.Range("H32").Value = 120
$ws.Range("I32").Value = 50
$ws.Range("J32").Value = 26.03447842597961
$ws.Range("K32").Value = 'BufferX'
$ws.Range("L32").Value = 40
$ws.Range("M32").Value = 1.826724137931035
$ws.Range("N32").Value = 'Lipid1'
$ws.Range("O32").Value = 'na'
$ws.Range("P32").Value = 40
$ws.Range("Q32").Value = 2.031878879310343
$ws.Range("R32").Value = 'Lipid2'
$ws.Range("S32").Value = 'na'
$ws.Range("T32").Value = 40
$ws.Range("U32").Value = 1.08920646551724
$ws.Range("V32").Value = 'Lipid3'
$ws.Range("W32").Value = 'na'
$ws.Range("X32").Value = 40
$ws.Range("Y32").Value = 1.239335775862069

# Row 33
$ws.Range("A33").Value = '0.3333333333333333-FlowTest150424nanana'
$ws.Range("B33").Value = 'Failed'
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = '10:13:13'
$ws.Range("AA1").Formula = "=""2024-04-16"""
$ws.Range("AA1").Copy()
$ws.Range("E33").PasteSpecial(-4163)
$ws.Range("F33").Value = '[1, 4]'
$ws.Range("G33").Value = 0.3333333333333333
$ws.Range("H33").Value = 120
$ws.Range("I33").Value = 50
$ws.Range("J33").Value = 31.4096953868866
$ws.Range("K33").Value = 'BufferX'
$ws.Range("L33").Value = 60
$ws.Range("M33").Value = 1.257471264367816
$ws.Range("N33").Value = 'Lipid1'
$ws.Range("O33").Value = 'na'
$ws.Range("P33").Value = 60
$ws.Range("Q33").Value = 0.9258215517241375
$ws.Range("R33").Value = 'Lipid2'
$ws.Range("S33").Value = 'na'
$ws.Range("T33").Value = 60
$ws.Range("U33").Value = 0.7070284482758619
$ws.Range("V33").Value = 'Lipid3'
$ws.Range("W33").Value = 'na'
$ws.Range("X33").Value = 60
$ws.Range("Y33").Value = 0.7396146551724132

# Row 34
$ws.Range("A34").Value = '0.3333333333333333-FlowTest150424nanana'
$ws.Range("B34").Value = 'Failed'
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = '10:14:01'
$ws.Range("AA1").Formula = "=""2024-04-16"""
$ws.Range("AA1").Copy()
$ws.Range("E34").PasteSpecial(-4163)
$ws.Range("F34").Value = '[1, 5]'
$ws.Range("G34").Value = 0.3333333333333333
$ws.Range("H34").Value = 120
$ws.Range("I34").Value = 50
$ws.Range("J34").Value = 30.26472544670105
$ws.Range("K34").Value = 'BufferX'
$ws.Range("L34").Value = 70
$ws.Range("M34").Value = 1.16551724137931
$ws.Range("N34").Value = 'Lipid1'
$ws.Range("O34").Value = 'na'
$ws.Range("P34").Value = 70
$ws.Range("Q34").Value = 0.7724036945812802
$ws.Range("R34").Value = 'Lipid2'
$ws.Range("S34").Value = 'na'
$ws.Range("T34").Value = 70
$ws.Range("U34").Value = 0.5429701970443349
$ws.Range("V34").Value = 'Lipid3'
$ws.Range("W34").Value = 'na'
$ws.Range("X34").Value = 70
$ws.Range("Y34").Value = 0.6088076354679803

# Row 35
$ws.Range("A35").Value = '0.3333333333333333-FlowTest150424nanana'
$ws.Range("B35").Value = 'Failed to Eq'
$ws.Range("C35").Value = 1
$ws.Range("D35").Value = '10:44:43'
$ws.Range("AA1").Formula = "=""2024-04-16"""
$ws.Range("AA1").Copy()
$ws.Range("E35").PasteSpecial(-4163)
$ws.Range("F35").Value = '[1, 2]'
$ws.Range("G35").Value = 0.3333333333333333
$ws.Range("H35").Value = 120
$ws.Range("I35").Value = 50
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 'BufferX'
$ws.Range("L35").Value = 50
$ws.Range("M35").Value = 2.520689655172414
$ws.Range("N35").Value = 'Lipid1'
$ws.Range("O35").Value = 'na'
$ws.Range("P35").Value = 50
$ws.Range("Q35").Value = 7.051565862068966
$ws.Range("R35").Value = 'Lipid2'
$ws.Range("S35").Value = 'na'
$ws.Range("T35").Value = 50
$ws.Range("U35").Value = 6.604669310344828
$ws.Range("V35").Value = 'Lipid3'
$ws.Range("W35").Value = 'na'
$ws.Range("X35").Value = 50
$ws.Range("Y35").Value = 6.694048620689656

# Row 36
$ws.Range("A36").Value = '0.3333333333333333-FlowTest150424nanana'
$ws.Range("B36").Value = 'Failed'
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = '10:50:40'
$ws.Range("AA1").Formula = "=""2024-04-16"""
$ws.Range("AA1").Copy()
$ws.Range("E36").PasteSpecial(-4163)
$ws.Range("F36").Value = '[1, 2]'
$ws.Range("G36").Value = 0.3333333333333333
$ws.Range("H36").Value = 120
$ws.Range("I36").Value = 50
$ws.Range("J36").Value = 32.52000260353088
$ws.Range("K36").Value = 'BufferX'
$ws.Range("L36").Value = 50
$ws.Range("M36").Value = 0.02896551724137936
$ws.Range("N36").Value = 'Lipid1'
$ws.Range("O36").Value = 'na'
$ws.Range("P36").Value = 50
$ws.Range("Q36").Value = 0.1183306896551724
$ws.Range("R36").Value = 'Lipid2'
$ws.Range("S36").Value = 'na'
$ws.Range("T36").Value = 50
$ws.Range("U36").Value = 0.101572068965517
$ws.Range("V36").Value = 'Lipid3'
$ws.Range("W36").Value = 'na'
$ws.Range("X36").Value = 50
$ws.Range("Y36").Value = 0.09039965517241327

# Row 37
$ws.Range("A37").Value = '0.3333333333333333-FlowTest150424nanana'
$ws.Range("B37").Value = 'Failed'
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = '10:56:36'
$ws.Range("AA1").Formula = "=""2024-04-16"""
$ws.Range("AA1").Copy()
$ws.Range("E37").PasteSpecial(-4163)
$ws.Range("F37").Value = '[1, 2]'
$ws.Range("G37").Value = 0.3333333333333333
$ws.Range("H37").Value = 120
$ws.Range("I37").Value = 50
$ws.Range("J37").Value = 37.01172184944153
$ws.Range("K37").Value = 'BufferX'
$ws.Range("L37").Value = 50
$ws.Range("M37").Value = 0.03310344827586206
$ws.Range("N37").Value = 'Lipid1'
$ws.Range("O37").Value = 'na'
$ws.Range("P37").Value = 50
$ws.Range("Q37").Value = 0.04087620689655211
$ws.Range("R37").Value = 'Lipid2'
$ws.Range("S37").Value = 'na'
$ws.Range("T37").Value = 50
$ws.Range("U37").Value = 0.05204862068965582
$ws.Range("V37").Value = 'Lipid3'
$ws.Range("W37").Value = 'na'
$ws.Range("X37").Value = 50
$ws.Range("Y37").Value = 0.04087620689655211

# Row 38
$ws.Range("A38").Value = '0.3333333333333333-FlowTest150424nanana'
$ws.Range("B38").Value = 'Failed'
$ws.Range("C38").Value = 1
$ws.Range("D38").Value = '10:57:35'
$ws.Range("AA1").Formula = "=""2024-04-16"""
$ws.Range("AA1").Copy()
$ws.Range("E38").PasteSpecial(-4163)
$ws.Range("F38").Value = '[1, 3]'
$ws.Range("G38").Value = 0.3333333333333333
$ws.Range("H38").Value = 120
$ws.Range("I38").Value = 50
$ws.Range("J38").Value = 28.42667031288147
$ws.Range("K38").Value = 'BufferX'
$ws.Range("L38").Value = 40
$ws.Range("M38").Value = 0.06465517241379307
$ws.Range("N38").Value = 'Lipid1'
$ws.Range("O38").Value = 'na'
$ws.Range("P38").Value = 40
$ws.Range("Q38").Value = 0.04877543103448261
$ws.Range("R38").Value = 'Lipid2'
$ws.Range("S38").Value = 'na'
$ws.Range("T38").Value = 40
$ws.Range("U38").Value = 0.05924956896551663
$ws.Range("V38").Value = 'Lipid3'
$ws.Range("W38").Value = 'na'
$ws.Range("X38").Value = 40
$ws.Range("Y38").Value = 0.05575818965517243

# Row 39
$ws.Range("A39").Value = '0.3333333333333333-FlowTest150424nanana'
$ws.Range("B39").Value = 'Failed'
$ws.Range("C39").Value = 1
$ws.Range("D39").Value = '10:58:24'
$ws.Range("AA1").Formula = "=""2024-04-16"""
$ws.Range("AA1").Copy()
$ws.Range("E39").PasteSpecial(-4163)
$ws.Range("F39").Value = '[1, 4]'
$ws.Range("G39").Value = 0.3333333333333333
$ws.Range("H39").Value = 120
$ws.Range("I39").Value = 50
$ws.Range("J39").Value = 28.78419852256775
$ws.Range("K39").Value = 'BufferX'
$ws.Range("L39").Value = 60
$ws.Range("M39").Value = 0.04712643678160925
$ws.Range("N39").Value = 'Lipid1'
$ws.Range("O39").Value = 'na'
$ws.Range("P39").Value = 60
$ws.Range("Q39").Value = 0.03779913793103518
$ws.Range("R39").Value = 'Lipid2'
$ws.Range("S39").Value = 'na'
$ws.Range("T39").Value = 60
$ws.Range("U39").Value = 0.03547155172413809
$ws.Range("V39").Value = 'Lipid3'
$ws.Range("W39").Value = 'na'
$ws.Range("X39").Value = 60
$ws.Range("Y39").Value = 0.05176465517241387

# Row 40
$ws.Range("A40").Value = '0.3333333333333333-FlowTest150424nanana'
$ws.Range("B40").Value = 'Failed'
$ws.Range("C40").Value = 1
$ws.Range("D40").Value = '10:59:04'
$ws.Range("AA1").Formula = "=""2024-04-16"""
$ws.Range("AA1").Copy()
$ws.Range("E40").PasteSpecial(-4163)
$ws.Range("F40").Value = '[1, 5]'
$ws.Range("G40").Value = 0.3333333333333333
$ws.Range("H40").Value = 120
$ws.Range("I40").Value = 50
$ws.Range("J40").Value = 21.44424200057983
$ws.Range("K40").Value = 'BufferX'
$ws.Range("L40").Value = 70
$ws.Range("M40").Value = 0.0261083743842365
$ws.Range("N40").Value = 'Lipid1'
$ws.Range("O40").Value = 'na'
$ws.Range("P40").Value = 70
$ws.Range("Q40").Value = 0.04956674876847294
$ws.Range("R40").Value = 'Lipid2'
$ws.Range("S40").Value = 'na'
$ws.Range("T40").Value = 70
$ws.Range("U40").Value = 0.03959137931034503
$ws.Range("V40").Value = 'Lipid3'
$ws.Range("W40").Value = 'na'
$ws.Range("X40").Value = 70
$ws.Range("Y40").Value = 0.03360615763546845

# Row 41
$ws.Range("A41").Value = '0.3333333333333333-FlowTest150424nanana'
$ws.Range("B41").Value = 'Failed'
$ws.Range("C41").Value = 1
$ws.Range("D41").Value = '10:59:47'
$ws.Range("AA1").Formula = "=""2024-04-16"""
$ws.Range("AA1").Copy()
$ws.Range("E41").PasteSpecial(-4163)
$ws.Range("F41").Value = '[1, 6]'
$ws.Range("G41").Value = 0.3333333333333333
$ws.Range("H41").Value = 120
$ws.Range("I41").Value = 50
$ws.Range("J41").Value = 28.24206376075745
$ws.Range("K41").Value = 'BufferX'
$ws.Range("L41").Value = 80
$ws.Range("M41").Value = 0.03232758620689662
$ws.Range("N41").Value = 'Lipid1'
$ws.Range("O41").Value = 'na'
$ws.Range("P41").Value = 80
$ws.Range("Q41").Value = 0.04268125000000023
$ws.Range("R41").Value = 'Lipid2'
$ws.Range("S41").Value = 'na'
$ws.Range("T41").Value = 80
$ws.Range("U41").Value = 0.0322071120689655
$ws.Range("V41").Value = 'Lipid3'
$ws.Range("W41").Value = 'na'
$ws.Range("X41").Value = 80
$ws.Range("Y41").Value = 0.0356984913793104

# Row 42
$ws.Range("A42").Value = '0.3333333333333333-FlowTest150424nanana'
$ws.Range("B42").Value = 'Failed'
$ws.Range("C42").Value = 1
$ws.Range("D42").Value = '11:01:11'
$ws.Range("AA1").Formula = "=""2024-04-16"""
$ws.Range("AA1").Copy()
$ws.Range("E42").PasteSpecial(-4163)
$ws.Range("F42").Value = '[1, 7]'
$ws.Range("G42").Value = 0.3333333333333333
$ws.Range("H42").Value = 120
$ws.Range("I42").Value = 50
$ws.Range("J42").Value = 43.78574204444885
$ws.Range("K42").Value = 'BufferX'
$ws.Range("L42").Value = 30
$ws.Range("M42").Value = 0.04022988505747129
$ws.Range("N42").Value = 'Lipid1'
$ws.Range("O42").Value = 'na'
$ws.Range("P42").Value = 30
$ws.Range("Q42").Value = 0.04387413793103541
$ws.Range("R42").Value = 'Lipid2'
$ws.Range("S42").Value = 'na'
$ws.Range("T42").Value = 30
$ws.Range("U42").Value = 0.1097465517241372
$ws.Range("V42").Value = 'Lipid3'
$ws.Range("W42").Value = 'na'
$ws.Range("X42").Value = 30
$ws.Range("Y42").Value = 0.1376775862068965

# Clean up helper cell used for date literals
$ws.Range("AA1").ClearContents()
